$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 32 with the next day's user impact data
$ws.Range("A32").Value = 45981
$ws.Range("B32").Value = 5609
$ws.Range("C32").Value = 4344
$ws.Range("D32").Value = 4035
$ws.Range("E32").Value = 228
$ws.Range("F32").Value = 43
$ws.Range("G32").Value = 29
$ws.Range("H32").Value = 8
$ws.Range("I32").Value = 1

# Match the date format style used by the preceding rows (row 31, column A)
$ws.Range("A32").NumberFormat = $ws.Range("A31").NumberFormat

# Move the active selection to the newly added row, matching prior behavior
$excel.Goto($ws.Range("A32:I32"))
